# Adding inverter to BOM
# Inverter was missing from BOM, corrected

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D8 was blank - set to "NA" to match corrected data
$ws.Range("D8").Value = "NA"

# Insert a new row at row 10 (pushes existing rows 10-19 down to 11-20).
$ws.Rows.Item(10).Insert()

# Row insert above pulls column default formats; re-apply the same cell
# styles used by the other data rows (copy format only from row 9, the
# row immediately above, restricted to the table's B:I columns).
$ws.Range("B9:I9").Copy()
$ws.Range("B10:I10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new Inverter row (row 10)
$ws.Cells.Item(10, 2).Value = "Inverter"                                    # B10 Part Name
$ws.Cells.Item(10, 3).Value = "Victron"                                     # C10 Manufacture
$ws.Cells.Item(10, 4).Value = "Phoenix 250VA"                               # D10 Manufacture Part Number
$ws.Cells.Item(10, 5).Value = "Phoenix 250VA 12-Volt 120V AC Pure Sine Wave Inverter"  # E10 Part Description
$ws.Cells.Item(10, 6).Value = 96.9                                          # F10 Unit Price
$ws.Cells.Item(10, 7).Value = 1                                             # G10 Quantity
$ws.Cells.Item(10, 8).Formula = "=F10*G10"                                  # H10 Extended Price
$ws.Cells.Item(10, 9).Value = "https://www.amazon.com/gp/product/B01NAO10QX/ref=ox_sc_act_image_1?smid=AERMGYAT5R869&th=1"  # I10 Supplier Link

# Add the hyperlink for the new supplier link cell
$ws.Hyperlinks.Add($ws.Cells.Item(10, 9), "https://www.amazon.com/gp/product/B01NAO10QX/ref=ox_sc_act_image_1?smid=AERMGYAT5R869&th=1") | Out-Null

# Hyperlinks.Add applies its own default style; restore the table's
# standard cell style (matching the other Supplier Link cells) on top.
$ws.Range("I9").Copy()
$ws.Range("I10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Selection as captured in the saved file
$ws.Range("I11").Select()

# Window position shifted slightly in the saved session
$excel.Width = $excel.Width
